$d = $word.ActiveDocument

# 1. "stimmst du zu" -> "nimmst du zur Kenntnis"
$d.Content.Find.Execute("stimmst du zu", $true, $false, $false, $false, $false,
                         $true, 1, $false, "nimmst du zur Kenntnis", 2)

# 2. "Du stimmst ebenfalls zu" -> "Du nimmst ebenfalls zur Kenntnis"
$d.Content.Find.Execute("Du stimmst ebenfalls zu", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Du nimmst ebenfalls zur Kenntnis", 2)

# 3. Uppercase "oder" -> "ODER" in the WEBSITE/ANWENDUNG disclaimer paragraph (first occurrence)
$d.Content.Find.Execute("DIE WEBSITE oder MOBILE ANWENDUNG ZUR VERFÜGUNG GESTELLT WERDEN", $true, $false, $false, $false, $false,
                         $true, 1, $false, "DIE WEBSITE ODER MOBILE ANWENDUNG ZUR VERFÜGUNG GESTELLT WERDEN", 2)

# 4. Uppercase "oder" -> "ODER" (second occurrence, "BETRIEB DER WEBSITE oder ANWENDUNG")
$d.Content.Find.Execute("BETRIEB DER WEBSITE oder ANWENDUNG ODER DER DIENSTE", $true, $false, $false, $false, $false,
                         $true, 1, $false, "BETRIEB DER WEBSITE ODER ANWENDUNG ODER DER DIENSTE", 2)

# 5. Remove "zu oder mobilen " from "zu ergänzen oder zu oder mobilen löschen"
$d.Content.Find.Execute("zu ergänzen oder zu oder mobilen löschen", $true, $false, $false, $false, $false,
                         $true, 1, $false, "zu ergänzen oder löschen", 2)
